# "Generate Report for Handback"
#
# The handback transform for file d7b95ec1-085b-40da-95a2-610af571cd68 (row 7
# in each status table) failed because the handback file name didn't match
# the handoff file name it was supposed to correspond to. Reflect that in the
# report:
#   - Overview sheet: update the zh-cn / de-de status cells for that file.
#   - zh-cn sheet: update its own Status cell and record the error detail.
#   - de-de sheet: update its own Status cell and record the error detail.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!B7 (zh-cn status) and Overview!C7 (de-de status) for the d7b95ec1 row
$wsOverview.Cells.Item(7, 2).Value = $newStatus
$wsOverview.Cells.Item(7, 3).Value = $newStatus

# zh-cn sheet, row 7 (d7b95ec1 file): Status column (C) + Error Detail column (L)
$wsZhCn.Cells.Item(7, 3).Value = $newStatus
$wsZhCn.Cells.Item(7, 12).Value = "Handback file name: kz3nv3g0.o1w is different with handoff file name: d7b95ec1-085b-40da-95a2-610af571cd68.657c1fa53195ba9de368c3b3f49c60d705ba7a44.zh-cn."

# de-de sheet, row 7 (d7b95ec1 file): Status column (C) + Error Detail column (L)
$wsDeDe.Cells.Item(7, 3).Value = $newStatus
$wsDeDe.Cells.Item(7, 12).Value = "Handback file name: kz3nv3g0.o1w is different with handoff file name: d7b95ec1-085b-40da-95a2-610af571cd68.657c1fa53195ba9de368c3b3f49c60d705ba7a44.de-de."
